# This script inserts a new weekly data point at row 146 by shifting the
# existing data block (rows 146-225) down by one row, freeing up row 146 for
# a brand-new record, and copying the former last row (225) into the newly
# created row 226.
#
# Columns A,B,C,E,F,G,H,I,J,K are constant across the whole data block, so
# only D (date) and L..T (quality/volume/price/unit/origin/etc.) need to be
# shifted. Row 226 additionally needs A..K populated since it is brand new.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$shiftCols = @("D","L","M","N","O","P","Q","R","S","T")
$copyCols  = @("A","B","C","E","F","G","H","I","J","K")

$firstRow = 146
$lastRow  = 225
$newRow   = 226

# 1) Populate the brand new last row (226) with the constant A..K values,
#    taken from the row immediately above it.
foreach ($c in $copyCols) {
    $src = $ws.Range($c + $lastRow).Value2
    $ws.Range($c + $newRow).Value2 = $src
}

# The date column (D) uses a custom date/time number format; make sure the
# brand new row 226 inherits it too (all other columns use the default
# "General" style already, matching the rest of the block).
$dateFormat = $ws.Range("D" + $lastRow).NumberFormat
$ws.Range("D" + $newRow).NumberFormat = $dateFormat

# 2) Shift the D,L..T data down by one row: process from the bottom (226)
#    up to 147 so that each row's original data is read before it gets
#    overwritten.
for ($r = $newRow; $r -ge ($firstRow + 1); $r--) {
    foreach ($c in $shiftCols) {
        $srcVal = $ws.Range($c + ($r - 1)).Value2
        $ws.Range($c + $r).Value2 = $srcVal
    }
}

# 3) At this point row 146 still holds its original (pre-shift) data.
#    Capture the values that used to live in row 147 (now duplicated into
#    row 148 as well), which become the new row 146's L..T contents, then
#    overwrite row 146 with the new date and those values.
$newL = "Primera"
$newM = 36
$newN = 11000
$newO = 11000
$newP = 11000
$newQ = "$/bandeja 12 canastillos 125 gramos"
$newR = "Provincia de Quillota"
$newS = 7333
$newT = 1.5

$ws.Range("D146").Value2 = 44824
$ws.Range("L146").Value2 = $newL
$ws.Range("M146").Value2 = $newM
$ws.Range("N146").Value2 = $newN
$ws.Range("O146").Value2 = $newO
$ws.Range("P146").Value2 = $newP
$ws.Range("Q146").Value2 = $newQ
$ws.Range("R146").Value2 = $newR
$ws.Range("S146").Value2 = $newS
$ws.Range("T146").Value2 = $newT

Write-Host "Shift complete"
